$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 40.66478
$ws.Range("B2").Value = -74.89136999999999
$ws.Range("A3").Value = 38.39657
$ws.Range("B3").Value = -98.70086000000001
$ws.Range("A4").Value = 40.34971
$ws.Range("B4").Value = -79.84963999999999
$ws.Range("A5").Value = 34.238356
$ws.Range("B5").Value = -118.253428
$ws.Range("A6").Value = 41.351626
$ws.Range("B6").Value = -85.13566
$ws.Range("A7").Value = 42.210608
$ws.Range("B7").Value = -93.215351
$ws.Range("A8").Value = 33.83575
$ws.Range("B8").Value = -118.11689
$ws.Range("A9").Value = 43.06623
$ws.Range("B9").Value = -92.68472
$ws.Range("A10").Value = 44.33037
$ws.Range("B10").Value = -75.90884
$ws.Range("A11").Value = 44.607389
$ws.Range("B11").Value = -93.78310500000001
$ws.Range("A12").Value = 44.5622
$ws.Range("B12").Value = -74.905897
$ws.Range("A13").Value = 34.00255
$ws.Range("B13").Value = -118.1327
$ws.Range("A14").Value = 42.84871
$ws.Range("B14").Value = -106.3146
$ws.Range("A15").Value = 39.281609
$ws.Range("B15").Value = -121.858902
$ws.Range("A16").Value = 39.1072
$ws.Range("B16").Value = -84.61387999999999
$ws.Range("A17").Value = 40.44519
$ws.Range("B17").Value = -74.93169
$ws.Range("A18").Value = 41.55228
$ws.Range("B18").Value = -73.07089000000001
$ws.Range("A19").Value = 32.65896
$ws.Range("B19").Value = -117.10381
$ws.Range("A20").Value = 41.18806
$ws.Range("B20").Value = -103.67739
$ws.Range("A21").Value = 36.3733
$ws.Range("B21").Value = -84.23925
$ws.Range("A22").Value = 41.765963
$ws.Range("B22").Value = -87.663831
$ws.Range("A23").Value = 38.737189
$ws.Range("B23").Value = -122.145336
$ws.Range("A24").Value = 34.320805
$ws.Range("B24").Value = -111.00948
$ws.Range("A25").Value = 41.028894
$ws.Range("B25").Value = -76.466869
$ws.Range("A26").Value = 40.55129
$ws.Range("B26").Value = -89.25015
$ws.Range("A27").Value = 45.56802
$ws.Range("B27").Value = -94.5924
$ws.Range("A28").Value = 45.54008
$ws.Range("B28").Value = -122.64731
$ws.Range("A29").Value = 34.224205
$ws.Range("B29").Value = -118.607902
$ws.Range("A30").Value = 44.54588
$ws.Range("B30").Value = -70.55643000000001
$ws.Range("A31").Value = 42.72435
$ws.Range("B31").Value = -92.19929
$ws.Range("A32").Value = 34.105383
$ws.Range("B32").Value = -117.950665
$ws.Range("A33").Value = 40.23321
$ws.Range("B33").Value = -74.70677000000001
$ws.Range("A34").Value = 46.90732
$ws.Range("B34").Value = -95.06321
$ws.Range("A35").Value = 40.74836
$ws.Range("B35").Value = -74.24467
$ws.Range("A36").Value = 38.39053
$ws.Range("B36").Value = -90.58468000000001
$ws.Range("A37").Value = 30.07394
$ws.Range("B37").Value = -92.68047
$ws.Range("A38").Value = 33.40379
$ws.Range("B38").Value = -82.10532000000001
$ws.Range("A39").Value = 39.56686
$ws.Range("B39").Value = -75.0137
$ws.Range("A40").Value = 33.228333
$ws.Range("B40").Value = -97.303611
$ws.Range("A41").Value = 40.635226
$ws.Range("B41").Value = -73.967628
$ws.Range("A42").Value = 42.676621
$ws.Range("B42").Value = -83.239902
$ws.Range("A43").Value = 41.88805
$ws.Range("B43").Value = -87.848208
$ws.Range("A44").Value = 25.91851
$ws.Range("B44").Value = -97.44678999999999
$ws.Range("A45").Value = 43.669225
$ws.Range("B45").Value = -92.972261
$ws.Range("A46").Value = 37.8088
$ws.Range("B46").Value = -101.58307
$ws.Range("A47").Value = 34.94638
$ws.Range("B47").Value = -81.03807
$ws.Range("A48").Value = 34.065556
$ws.Range("B48").Value = -117.047222
$ws.Range("A49").Value = 41.750688
$ws.Range("B49").Value = -87.68286500000001
$ws.Range("A50").Value = 36.61676
$ws.Range("B50").Value = -88.32096
$ws.Range("A51").Value = 40.17353
$ws.Range("B51").Value = -74.02208
$ws.Range("A52").Value = 38.54863
$ws.Range("B52").Value = -90.37566
$ws.Range("A53").Value = 41.059085
$ws.Range("B53").Value = -73.765974
$ws.Range("A54").Value = 44.8898
$ws.Range("B54").Value = -89.51094000000001
$ws.Range("A55").Value = 43.3463
$ws.Range("B55").Value = -83.384
$ws.Range("A56").Value = 30.182358
$ws.Range("B56").Value = -85.72841699999999
$ws.Range("A57").Value = 37.506654
$ws.Range("B57").Value = -122.252962
$ws.Range("A58").Value = 45.51544
$ws.Range("B58").Value = -122.65868
$ws.Range("A59").Value = 29.467922
$ws.Range("B59").Value = -95.11463000000001
$ws.Range("A60").Value = 36.911944
$ws.Range("B60").Value = -83.919444
$ws.Range("A61").Value = 41.16786
$ws.Range("B61").Value = -111.96861
$ws.Range("A62").Value = 41.56639
$ws.Range("B62").Value = -73.918181
$ws.Range("A63").Value = 38.57299
$ws.Range("B63").Value = -121.49938
$ws.Range("A64").Value = 39.188043
$ws.Range("B64").Value = -120.834568
$ws.Range("A65").Value = 40.92497
$ws.Range("B65").Value = -72.69535999999999
$ws.Range("A66").Value = 33.368661
$ws.Range("B66").Value = -87.01028700000001
$ws.Range("A67").Value = 40.88991
$ws.Range("B67").Value = -124.08757
$ws.Range("A68").Value = 46.666278
$ws.Range("B68").Value = -112.514167
$ws.Range("A69").Value = 39.478111
$ws.Range("B69").Value = -96.07643899999999
$ws.Range("A70").Value = 33.699624
$ws.Range("B70").Value = -78.910265
$ws.Range("A71").Value = 42.43759
$ws.Range("B71").Value = -77.096585
$ws.Range("A72").Value = 33.75575
$ws.Range("B72").Value = -117.87106
$ws.Range("A73").Value = 40.647771
$ws.Range("B73").Value = -79.107394
$ws.Range("A74").Value = 44.865604
$ws.Range("B74").Value = -93.39997
$ws.Range("A75").Value = 41.04158
$ws.Range("B75").Value = -96.10195
$ws.Range("A76").Value = 40.7899
$ws.Range("B76").Value = -74.02005
$ws.Range("A77").Value = 40.73868
$ws.Range("B77").Value = -73.98324
$ws.Range("A78").Value = 47.70375
$ws.Range("B78").Value = -116.785747
$ws.Range("A79").Value = 39.06548
$ws.Range("B79").Value = -84.11596
$ws.Range("A80").Value = 38.59822
$ws.Range("B80").Value = -90.16107
$ws.Range("A81").Value = 26.68223
$ws.Range("B81").Value = -80.11203999999999
$ws.Range("A82").Value = 38.19563
$ws.Range("B82").Value = -86.97767
$ws.Range("A83").Value = 37.35819
$ws.Range("B83").Value = -120.74194
$ws.Range("A84").Value = 33.5096
$ws.Range("B84").Value = -112.1026
$ws.Range("A85").Value = 43.64384
$ws.Range("B85").Value = -70.98053
$ws.Range("A86").Value = 42.45838
$ws.Range("B86").Value = -93.80831000000001
$ws.Range("A87").Value = 30.362953
$ws.Range("B87").Value = -87.13959199999999
$ws.Range("A88").Value = 39.62418
$ws.Range("B88").Value = -86.47978999999999
$ws.Range("A89").Value = 41.85287
$ws.Range("B89").Value = -71.39655
$ws.Range("A90").Value = 33.727633
$ws.Range("B90").Value = -118.066351
$ws.Range("A91").Value = 40.1942
$ws.Range("B91").Value = -80.3128
$ws.Range("A92").Value = 33.18422
$ws.Range("B92").Value = -117.29259
$ws.Range("A93").Value = 40.583031
$ws.Range("B93").Value = -79.76570700000001
$ws.Range("A94").Value = 45.56393
$ws.Range("B94").Value = -94.94584999999999
$ws.Range("A95").Value = 33.72709
$ws.Range("B95").Value = -116.39837
$ws.Range("A96").Value = 42.727046
$ws.Range("B96").Value = -84.555521
$ws.Range("A97").Value = 39.95224
$ws.Range("B97").Value = -74.985556
$ws.Range("A98").Value = 37.52802
$ws.Range("B98").Value = -122.02633
$ws.Range("A99").Value = 40.33382
$ws.Range("B99").Value = -74.04594
$ws.Range("A100").Value = 36.8879
$ws.Range("B100").Value = -111.44425
$ws.Range("A101").Value = 41.82161
$ws.Range("B101").Value = -71.35842
